$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Email" column (header + value)
$ws.Range("H1").Value = "Email"
$ws.Range("H2").Value = "adrianrentea01@gmail.com"

# Make the new column fit its content, like the other bestFit columns
# (23.666... round-trips through the engine's internal width conversion
# to the target displayed column width of 24.5)
$ws.Columns.Item(8).ColumnWidth = 23.666666666666668

# Update the selected cell on the sheet to match the new edit location
$ws.Range("H8").Select() | Out-Null
